$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 846.92
$ws.Range("J17").Value = 869.8333
$ws.Range("L17").Value = 2609.4999
$ws.Range("N17").Value = -2945.4999
$ws.Range("H137").Value = 1409.9722
$ws.Range("I137").Value = 1192.3214
$ws.Range("J137").Value = 2171.75
$ws.Range("K137").Value = 3576.9642
$ws.Range("L137").Value = 6515.25
$ws.Range("M137").Value = -1026.9642
$ws.Range("N137").Value = -11615.25

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24918.85
$ws.Range("I32").Value = 4143.645
$ws.Range("J32").Value = 142015.45
$ws.Range("K32").Value = 4143.645
$ws.Range("L32").Value = 142015.45
$ws.Range("M32").Value = -3856.645
$ws.Range("N32").Value = -142589.45
$ws.Range("H74").Value = 5002226.5
$ws.Range("I74").Value = 1611.1666
$ws.Range("J74").Value = 12503150
$ws.Range("K74").Value = 1611.1666
$ws.Range("L74").Value = 12503150
$ws.Range("M74").Value = -737.1666
$ws.Range("N74").Value = -12504898
$ws.Range("H77").Value = 5002226.5
$ws.Range("I77").Value = 1611.1666
$ws.Range("J77").Value = 12503150
$ws.Range("K77").Value = 8055.833000000001
$ws.Range("L77").Value = 62515750
$ws.Range("M77").Value = -3687.833000000001
$ws.Range("N77").Value = -62524486
$ws.Range("H80").Value = 18403.715
$ws.Range("J80").Value = 23765.2
$ws.Range("L80").Value = 23765.2
$ws.Range("N80").Value = -25761.2
$ws.Range("H83").Value = 18403.715
$ws.Range("J83").Value = 23765.2
$ws.Range("L83").Value = 71295.60000000001
$ws.Range("N83").Value = -81279.60000000001
$ws.Range("H86").Value = 30133
$ws.Range("I86").Value = 10285
$ws.Range("J86").Value = 49981
$ws.Range("K86").Value = 10285
$ws.Range("L86").Value = 49981
$ws.Range("M86").Value = -9099
$ws.Range("N86").Value = -52353
$ws.Range("H89").Value = 30133
$ws.Range("I89").Value = 10285
$ws.Range("J89").Value = 49981
$ws.Range("K89").Value = 30855
$ws.Range("L89").Value = 149943
$ws.Range("M89").Value = -24927
$ws.Range("N89").Value = -161799
$ws.Range("H122").Value = 1169.1154
$ws.Range("I122").Value = 1087.6957
$ws.Range("J122").Value = 1793.3334
$ws.Range("K122").Value = 3263.0871
$ws.Range("L122").Value = 5380.0002
$ws.Range("M122").Value = -813.0870999999997
$ws.Range("N122").Value = -10280.0002

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 34972.5
$ws.Range("J70").Value = 34972.5
$ws.Range("L70").Value = 34972.5
$ws.Range("N70").Value = -35558.5
$ws.Range("H73").Value = 34972.5
$ws.Range("J73").Value = 34972.5
$ws.Range("L73").Value = 34972.5
$ws.Range("N73").Value = -37000.5
$ws.Range("H94").Value = 859.8333
$ws.Range("I94").Value = 789.75
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 789.75
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -338.75
$ws.Range("N94").Value = -1902
$ws.Range("H99").Value = 1737.2727
$ws.Range("I99").Value = 1642.8572
$ws.Range("K99").Value = 1642.8572
$ws.Range("M99").Value = -144.8571999999999
$ws.Range("H105").Value = 88890
$ws.Range("I105").Value = 101558.8
$ws.Range("J105").Value = 79144.766
$ws.Range("K105").Value = 101558.8
$ws.Range("L105").Value = 79144.766
$ws.Range("M105").Value = -99811.8
$ws.Range("N105").Value = -82638.766

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 516.7
$ws.Range("I22").Value = 370.375
$ws.Range("J22").Value = 1102
$ws.Range("K22").Value = 370.375
$ws.Range("L22").Value = 1102
$ws.Range("M22").Value = -20.375
$ws.Range("N22").Value = -1802
$ws.Range("H31").Value = 33849.246
$ws.Range("I31").Value = 1159.2759
$ws.Range("K31").Value = 1159.2759
$ws.Range("M31").Value = -864.2759000000001
$ws.Range("H34").Value = 33849.246
$ws.Range("I34").Value = 1159.2759
$ws.Range("K34").Value = 1159.2759
$ws.Range("M34").Value = -957.2759000000001
$ws.Range("H74").Value = 28618.8
$ws.Range("J74").Value = 28618.8
$ws.Range("L74").Value = 28618.8
$ws.Range("N74").Value = -30366.8
$ws.Range("H77").Value = 28618.8
$ws.Range("J77").Value = 28618.8
$ws.Range("L77").Value = 85856.39999999999
$ws.Range("N77").Value = -94592.39999999999
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H93").Value = 16549
$ws.Range("I93").Value = 4000
$ws.Range("J93").Value = 24915
$ws.Range("K93").Value = 4000
$ws.Range("L93").Value = 24915
$ws.Range("M93").Value = -2128
$ws.Range("N93").Value = -28659
$ws.Range("H96").Value = 20124
$ws.Range("J96").Value = 20124
$ws.Range("L96").Value = 20124
$ws.Range("N96").Value = -25616

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 24921.334
$ws.Range("J37").Value = 24921.334
$ws.Range("L37").Value = 74764.00199999999
$ws.Range("N37").Value = -74988.00199999999
$ws.Range("H51").Value = 3242.7144
$ws.Range("J51").Value = 3242.7144
$ws.Range("L51").Value = 9728.143199999999
$ws.Range("N51").Value = -10648.1432
$ws.Range("H122").Value = 575.4375
$ws.Range("I122").Value = 456.75
$ws.Range("J122").Value = 694.125
$ws.Range("K122").Value = 4110.75
$ws.Range("L122").Value = 6247.125
$ws.Range("M122").Value = -1660.75
$ws.Range("N122").Value = -11147.125
$ws.Range("H131").Value = 774.58
$ws.Range("J131").Value = 831.1279
$ws.Range("L131").Value = 2493.3837
$ws.Range("N131").Value = -12573.3837

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 11226.667
$ws.Range("J52").Value = 11226.667
$ws.Range("L52").Value = 11226.667
$ws.Range("N52").Value = -11744.667
$ws.Range("H55").Value = 4378.364
$ws.Range("I55").Value = 2024
$ws.Range("J55").Value = 5723.7144
$ws.Range("K55").Value = 2024
$ws.Range("L55").Value = 5723.7144
$ws.Range("M55").Value = -1697
$ws.Range("N55").Value = -6377.7144

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3938000.5
$ws.Range("I16").Value = 5250290.5
$ws.Range("J16").Value = 1130.25
$ws.Range("K16").Value = 5250290.5
$ws.Range("L16").Value = 1130.25
$ws.Range("M16").Value = -5250120.5
$ws.Range("N16").Value = -1470.25
$ws.Range("H45").Value = 6456.32
$ws.Range("I45").Value = 4832.6
$ws.Range("K45").Value = 4832.6
$ws.Range("M45").Value = -4425.6
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 15336.75
$ws.Range("I74").Value = 14742
$ws.Range("J74").Value = 19500
$ws.Range("K74").Value = 14742
$ws.Range("L74").Value = 19500
$ws.Range("M74").Value = -13744
$ws.Range("N74").Value = -21496
$ws.Range("H77").Value = 15336.75
$ws.Range("I77").Value = 14742
$ws.Range("J77").Value = 19500
$ws.Range("K77").Value = 44226
$ws.Range("L77").Value = 58500
$ws.Range("M77").Value = -39234
$ws.Range("N77").Value = -68484
$ws.Range("H82").Value = 1131.4445
$ws.Range("I82").Value = 718.7143
$ws.Range("J82").Value = 1394.091
$ws.Range("K82").Value = 718.7143
$ws.Range("L82").Value = 1394.091
$ws.Range("M82").Value = -357.7143
$ws.Range("N82").Value = -2116.091
$ws.Range("H85").Value = 1131.4445
$ws.Range("I85").Value = 718.7143
$ws.Range("J85").Value = 1394.091
$ws.Range("K85").Value = 718.7143
$ws.Range("L85").Value = 1394.091
$ws.Range("M85").Value = 529.2857
$ws.Range("N85").Value = -3890.091
$ws.Range("H96").Value = 14248.833
$ws.Range("J96").Value = 14248.833
$ws.Range("L96").Value = 14248.833
$ws.Range("N96").Value = -19740.833

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 250700
$ws.Range("I100").Value = 333933.34
$ws.Range("K100").Value = 667866.6800000001
$ws.Range("M100").Value = -667325.6800000001

Write-Output "Aegis_Profits market-data refresh applied across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR."
